# aggiornamento 15, 16, 17 marzo
# Append three new daily rows (227-229) to the COVID tracking sheet,
# continuing the date series (Excel serials 44301, 44302, 44303) with
# their "nuovi pos." / "somma mobile 7gg." / "somma mobile 7gg. per
# 100mila abitanti" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 226 (the last existing data row) carries the date number format
# (style index 2) on column A. Copy that formatting down into the new
# rows so the new dates render/serialize the same way as the rest of
# the column.
$ws.Range("A226").Copy()
$ws.Range("A227:A229").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @{ Row = 227; Date = 44301; B = 1; C = 1; D = 33.71544167228591 },
    @{ Row = 228; Date = 44302; B = 0; C = 1; D = 33.71544167228591 },
    @{ Row = 229; Date = 44303; B = 0; C = 1; D = 33.71544167228591 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
